$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 21:50"

# --- Reorder countries (ranking changes) ---
$ws.Range("A52").Value = "Colombia"
$ws.Range("A53").Value = "Serbia"
$ws.Range("A54").Value = "Croacia"
$ws.Range("A55").Value = "Eslovenia"
$ws.Range("A65").Value = "Ucrania"
$ws.Range("A66").Value = "Marruecos"
$ws.Range("A67").Value = "Barein"
$ws.Range("A79").Value = "Costa Rica"
$ws.Range("A80").Value = "Kazajistan"
$ws.Range("A135").Value = "Guatemala"
$ws.Range("A136").Value = "Polinesia Francesa"
$ws.Range("A137").Value = "Jamaica"

# --- Updated case numbers ---
# Row 12
$ws.Range("B12").Value = 16605
$ws.Range("C12").Value = 683
$ws.Range("E12").Value = 14349
$ws.Range("G12").Value = 74
$ws.Range("H12").Value = 433

# Row 21
$ws.Range("B21").Value = 4830
$ws.Range("C21").Value = 200
$ws.Range("E21").Value = 4528
$ws.Range("G21").Value = 12
$ws.Range("H21").Value = 175

# Row 22
$ws.Range("B22").Value = 4635
$ws.Range("C22").Value = 190
$ws.Range("E22").Value = 4583

# Row 49
$ws.Range("E49").Value = 641
$ws.Range("G49").Value = 6
$ws.Range("H49").Value = 30

# Row 52
$ws.Range("B52").Value = 906
$ws.Range("C52").Value = 108
$ws.Range("D52").Value = 31
$ws.Range("E52").Value = 859
$ws.Range("F52").Value = 29
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 16

# Row 53
$ws.Range("B53").Value = 900
$ws.Range("C53").Value = 115
$ws.Range("D53").Value = 42
$ws.Range("E53").Value = 835
$ws.Range("F53").Value = 62
$ws.Range("G53").Value = 7
$ws.Range("H53").Value = 23

# Row 54
$ws.Range("B54").Value = 867
$ws.Range("C54").Value = 77
$ws.Range("D54").Value = 67
$ws.Range("E54").Value = 794
$ws.Range("F54").Value = 32
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 6

# Row 55
$ws.Range("B55").Value = 802
$ws.Range("C55").Value = 46
$ws.Range("D55").Value = 10
$ws.Range("E55").Value = 777
$ws.Range("F55").Value = 24
$ws.Range("G55").Value = 4
$ws.Range("H55").Value = 15

# Row 65
$ws.Range("B65").Value = 618
$ws.Range("C65").Value = 70
$ws.Range("D65").Value = 8
$ws.Range("E65").Value = 597
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 13

# Row 66
$ws.Range("B66").Value = 602
$ws.Range("C66").Value = 46
$ws.Range("D66").Value = 24
$ws.Range("E66").Value = 542
$ws.Range("F66").Value = 1
$ws.Range("G66").Value = 3
$ws.Range("H66").Value = 36

# Row 67
$ws.Range("B67").Value = 567
$ws.Range("C67").Value = 52
$ws.Range("D67").Value = 295
$ws.Range("E67").Value = 268
$ws.Range("F67").Value = 2
$ws.Range("H67").Value = 4

# Row 79
$ws.Range("B79").Value = 347
$ws.Range("C79").Value = 17
$ws.Range("D79").Value = 4
$ws.Range("E79").Value = 341
$ws.Range("F79").Value = 8
$ws.Range("G79").Value = 0

# Row 80
$ws.Range("B80").Value = 340
$ws.Range("C80").Value = 38
$ws.Range("D80").Value = 22
$ws.Range("E80").Value = 316
$ws.Range("F80").Value = 6
$ws.Range("G80").Value = 1

# Row 135
$ws.Range("B135").Value = 38
$ws.Range("C135").Value = 2
$ws.Range("D135").Value = 10
$ws.Range("E135").Value = 27
$ws.Range("H135").Value = 1

# Row 136
$ws.Range("D136").Value = 0
$ws.Range("E136").Value = 36
$ws.Range("F136").Value = 1
$ws.Range("H136").Value = 0

# Row 137
$ws.Range("D137").Value = 10
$ws.Range("E137").Value = 33
$ws.Range("F137").Value = 0
